# Scheduled-runner refresh of cached market-board / profit figures.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and the derived
# LeveProfit(NQ/HQ) columns (H:N) for a handful of leves across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching a fresh Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 276.35483
$ws.Range("I19").Value = 240.14285
$ws.Range("J19").Value = 306.17648
$ws.Range("K19").Value = 240.14285
$ws.Range("L19").Value = 306.17648
$ws.Range("M19").Value = -65.14285000000001
$ws.Range("N19").Value = -656.1764800000001
$ws.Range("H29").Value = 1800
$ws.Range("J29").Value = 2780
$ws.Range("L29").Value = 8340
$ws.Range("N29").Value = -8902
$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 6000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 6000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -7872
$ws.Range("H88").Value = 726757.7
$ws.Range("I88").Value = 487
$ws.Range("J88").Value = 1372331.6
$ws.Range("K88").Value = 487
$ws.Range("L88").Value = 1372331.6
$ws.Range("M88").Value = -81
$ws.Range("N88").Value = -1373143.6
$ws.Range("H91").Value = 726757.7
$ws.Range("I91").Value = 487
$ws.Range("J91").Value = 1372331.6
$ws.Range("K91").Value = 487
$ws.Range("L91").Value = 1372331.6
$ws.Range("M91").Value = 917
$ws.Range("N91").Value = -1375139.6
$ws.Range("H98").Value = 3243.4285
$ws.Range("I98").Value = 3550.3333
$ws.Range("J98").Value = 1402
$ws.Range("K98").Value = 3550.3333
$ws.Range("L98").Value = 1402
$ws.Range("M98").Value = -2052.3333
$ws.Range("N98").Value = -4398
$ws.Range("H122").Value = 3243.4285
$ws.Range("I122").Value = 3550.3333
$ws.Range("J122").Value = 1402
$ws.Range("K122").Value = 10650.9999
$ws.Range("L122").Value = 4206
$ws.Range("M122").Value = -8200.999899999999
$ws.Range("N122").Value = -9106
$ws.Range("H125").Value = 2050.7646
$ws.Range("I125").Value = 603
$ws.Range("J125").Value = 3679.5
$ws.Range("K125").Value = 5427
$ws.Range("L125").Value = 33115.5
$ws.Range("M125").Value = -2967
$ws.Range("N125").Value = -38035.5
$ws.Range("H132").Value = 11503552
$ws.Range("I132").Value = 16675455
$ws.Range("J132").Value = 10435.333
$ws.Range("K132").Value = 50026365
$ws.Range("L132").Value = 31305.999
$ws.Range("M132").Value = -50023835
$ws.Range("N132").Value = -36365.999
$ws.Range("H137").Value = 1612.119
$ws.Range("I137").Value = 1488.6786
$ws.Range("J137").Value = 1859
$ws.Range("K137").Value = 4466.0358
$ws.Range("L137").Value = 5577
$ws.Range("M137").Value = -1916.0358
$ws.Range("N137").Value = -10677
$ws.Range("H138").Value = 450048.9
$ws.Range("I138").Value = 1663.8422
$ws.Range("J138").Value = 559270.9
$ws.Range("K138").Value = 4991.5266
$ws.Range("L138").Value = 1677812.7
$ws.Range("M138").Value = 148.4733999999999
$ws.Range("N138").Value = -1688092.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6913.4824
$ws.Range("I32").Value = 5542.463
$ws.Range("K32").Value = 5542.463
$ws.Range("M32").Value = -5255.463
$ws.Range("H45").Value = 1211.5
$ws.Range("I45").Value = 1156.1
$ws.Range("J45").Value = 1350
$ws.Range("K45").Value = 1156.1
$ws.Range("L45").Value = 1350
$ws.Range("M45").Value = -779.0999999999999
$ws.Range("N45").Value = -2104
$ws.Range("H61").Value = 55556776
$ws.Range("I61").Value = 62501064
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 62501064
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -62500852
$ws.Range("N61").Value = -2924
$ws.Range("H122").Value = 2204.4119
$ws.Range("I122").Value = 1811.3846
$ws.Range("K122").Value = 5434.1538
$ws.Range("M122").Value = -2984.1538
$ws.Range("H136").Value = 55556776
$ws.Range("I136").Value = 62501064
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 187503192
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -187500642
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 350
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -27
$ws.Range("N22").Value = -746
$ws.Range("H64").Value = 419.6
$ws.Range("I64").Value = 413.7143
$ws.Range("K64").Value = 413.7143
$ws.Range("M64").Value = -188.7143
$ws.Range("H67").Value = 419.6
$ws.Range("I67").Value = 413.7143
$ws.Range("K67").Value = 413.7143
$ws.Range("M67").Value = 366.2857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1545.8043
$ws.Range("I31").Value = 1439.1163
$ws.Range("J31").Value = 3075
$ws.Range("K31").Value = 1439.1163
$ws.Range("L31").Value = 3075
$ws.Range("M31").Value = -1144.1163
$ws.Range("N31").Value = -3665
$ws.Range("H34").Value = 1545.8043
$ws.Range("I34").Value = 1439.1163
$ws.Range("J34").Value = 3075
$ws.Range("K34").Value = 1439.1163
$ws.Range("L34").Value = 3075
$ws.Range("M34").Value = -1237.1163
$ws.Range("N34").Value = -3479
$ws.Range("H108").Value = 28681
$ws.Range("J108").Value = 34261.332
$ws.Range("L108").Value = 34261.332
$ws.Range("N108").Value = -41941.332
$ws.Range("H122").Value = 1349.75
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
$ws.Range("H134").Value = 17243112
$ws.Range("I134").Value = 1601.2632
$ws.Range("J134").Value = 50001984
$ws.Range("K134").Value = 4803.7896
$ws.Range("L134").Value = 150005952
$ws.Range("M134").Value = -2268.7896
$ws.Range("N134").Value = -150011022

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3900228
$ws.Range("I4").Value = 2500093.5
$ws.Range("J4").Value = 6700497
$ws.Range("K4").Value = 7500280.5
$ws.Range("L4").Value = 20101491
$ws.Range("M4").Value = -7500168.5
$ws.Range("N4").Value = -20101715
$ws.Range("H113").Value = 702.5
$ws.Range("I113").Value = 587.875
$ws.Range("J113").Value = 748.35
$ws.Range("K113").Value = 1763.625
$ws.Range("L113").Value = 2245.05
$ws.Range("M113").Value = 406.375
$ws.Range("N113").Value = -6585.05
$ws.Range("H122").Value = 1856.8667
$ws.Range("I122").Value = 849.6667
$ws.Range("J122").Value = 2108.6667
$ws.Range("K122").Value = 7647.0003
$ws.Range("L122").Value = 18978.0003
$ws.Range("M122").Value = -5197.0003
$ws.Range("N122").Value = -23878.0003
$ws.Range("H131").Value = 27030452
$ws.Range("I131").Value = 76923540
$ws.Range("J131").Value = 5029.5835
$ws.Range("K131").Value = 230770620
$ws.Range("L131").Value = 15088.7505
$ws.Range("M131").Value = -230765580
$ws.Range("N131").Value = -25168.7505
$ws.Range("H139").Value = 2105.919
$ws.Range("I139").Value = 2463.7778
$ws.Range("J139").Value = 1766.8948
$ws.Range("K139").Value = 7391.3334
$ws.Range("L139").Value = 5300.6844
$ws.Range("M139").Value = -2251.3334
$ws.Range("N139").Value = -15580.6844
$ws.Range("H140").Value = 3111.587
$ws.Range("I140").Value = 2643.5
$ws.Range("K140").Value = 7930.5
$ws.Range("M140").Value = -2750.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1944.35
$ws.Range("I102").Value = 1878.2
$ws.Range("J102").Value = 2142.8
$ws.Range("K102").Value = 1878.2
$ws.Range("L102").Value = 2142.8
$ws.Range("M102").Value = -256.2
$ws.Range("N102").Value = -5386.8
$ws.Range("H126").Value = 1941.3889
$ws.Range("I126").Value = 1694.5555
$ws.Range("J126").Value = 2188.2222
$ws.Range("K126").Value = 5083.666499999999
$ws.Range("L126").Value = 6564.6666
$ws.Range("M126").Value = -2613.666499999999
$ws.Range("N126").Value = -11504.6666
$ws.Range("H134").Value = 26745.857
$ws.Range("I134").Value = 10000
$ws.Range("J134").Value = 29536.834
$ws.Range("K134").Value = 30000
$ws.Range("L134").Value = 88610.50199999999
$ws.Range("M134").Value = -27465
$ws.Range("N134").Value = -93680.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6279.273
$ws.Range("I40").Value = 3069
$ws.Range("J40").Value = 7483.125
$ws.Range("K40").Value = 3069
$ws.Range("L40").Value = 7483.125
$ws.Range("M40").Value = -2933
$ws.Range("N40").Value = -7755.125
$ws.Range("H132").Value = 2520.276
$ws.Range("I132").Value = 2283.5
$ws.Range("J132").Value = 2907.7273
$ws.Range("K132").Value = 6850.5
$ws.Range("L132").Value = 8723.1819
$ws.Range("M132").Value = -4320.5
$ws.Range("N132").Value = -13783.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 473.26666
$ws.Range("I113").Value = 367.44446
$ws.Range("J113").Value = 632
$ws.Range("K113").Value = 1102.33338
$ws.Range("L113").Value = 1896
$ws.Range("M113").Value = 1067.66662
$ws.Range("N113").Value = -6236
$ws.Range("H126").Value = 43479196
$ws.Range("I126").Value = 62500468
$ws.Range("J126").Value = 2002
$ws.Range("K126").Value = 187501404
$ws.Range("L126").Value = 6006
$ws.Range("M126").Value = -187498934
$ws.Range("N126").Value = -10946
$ws.Range("H133").Value = 34971.668
$ws.Range("J133").Value = 34971.668
$ws.Range("L133").Value = 34971.668
$ws.Range("N133").Value = -45091.668
$ws.Range("H136").Value = 1603.9333
$ws.Range("I136").Value = 1326.4
$ws.Range("J136").Value = 2159
$ws.Range("K136").Value = 3979.2
$ws.Range("L136").Value = 6477
$ws.Range("M136").Value = -1429.2
$ws.Range("N136").Value = -11577
